# Apply cryptos list update (Thu Aug 10 18:39:40 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# All target cells hold plain text (coin names, links, price & volume strings),
# so force Text number format before assigning to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.0000") into numbers and losing formatting.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '29.402.44'

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -0.37%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.847.40'

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9989'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '240.86'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.95%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.6327'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -3.22%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.0000'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.06%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.07593'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +1.08%  '

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -0.37%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '24.50'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '2.460.48'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +32.81%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.07722'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +1.17%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.565.15'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +21.09%  '

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -0.80%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.6862'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '82.99'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -0.90%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.000009909'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +4.29%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '6.172'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.69%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '29.438.33'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -0.35%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '231.59'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -2.49%  '

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -0.66%  '

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '7.608'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -1.19%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.000'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '154.36'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -1.66%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.1392'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -2.01%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.462'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -0.50%  '

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -0.78%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.474'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -0.85%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.05815'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -3.79%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.258'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +0.45%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.129'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -0.22%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.027'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -1.21%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.872'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.75%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.160'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -1.84%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.7180'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.76%  '

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.517.93'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +24.20%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.247.94'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +3.80%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.791'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -0.41%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.01809'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +1.39%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.9077'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.16%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '6.125'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -1.99%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.9992'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '

$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'Quant'

$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '101.48'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -0.47%  '

$c = $ws.Range('B46')
$c.NumberFormat = '@'
$c.Value = 'Aave'

$c = $ws.Range('C46')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '67.28'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +1.16%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '7.292'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -1.62%  '

$c = $ws.Range('B48')
$c.NumberFormat = '@'
$c.Value = 'EnergySwap'

$c = $ws.Range('C48')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '9.196'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +1.18%  '

$c = $ws.Range('B49')
$c.NumberFormat = '@'
$c.Value = 'TheSandbox'

$c = $ws.Range('C49')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.4015'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -1.01%  '

$c = $ws.Range('B50')
$c.NumberFormat = '@'
$c.Value = 'RenderToken'

$c = $ws.Range('C50')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.695'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +1.86%  '

$c = $ws.Range('B51')
$c.NumberFormat = '@'
$c.Value = 'Algorand'

$c = $ws.Range('C51')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.1123'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -0.39%  '
